$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.691.16"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.848.22"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.95"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.10%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4277"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -6.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3651"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.07%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.88"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.72%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07242"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.57%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8996"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -7.83%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.69"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -6.96%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.845.81"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.68%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.378"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.64%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.579"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.92%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06847"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "77.84"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -7.75%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008954"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.75%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9997"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -7.32%  "
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.668.91"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.17%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.976"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.91%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.65"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.22%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.069.57"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.02%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.049"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.59"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.32"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.03%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.352"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.82%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.53"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("B31").Value = "LidoDAOToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.799"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08953"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7877"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -9.54%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.528"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -11.25%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.927"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.38%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.093"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -12.23%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.000"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05459"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.099"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.46%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.971"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01930"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.39%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5066"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -7.99%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.812"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.31%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1640"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.76%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.315"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -11.12%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06634"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.70"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.40"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.35%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4717"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -8.64%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9999"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.647"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.72%  "
